$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value for the "Price" column (D) which is stored
# as text in the workbook. Force NumberFormat to Text ("@") before writing
# so numeric-looking strings are not reinterpreted as numbers by Excel.
$priceUpdates = @{
    "D2"  = "245.87"
    "D3"  = "22.09"
    "D4"  = "5.367"
    "D5"  = "0.05856"
    "D6"  = "3.394"
    "D7"  = "6.360"
    "D8"  = "0.8134"
    "D9"  = "1.026"
    "D10" = "0.1421"
    "D11" = "0.04004"
    "D12" = "0.07369"
    "D13" = "0.03022"
    "D14" = "4.162"
    "D15" = "0.09404"
    "D16" = "0.001590"
    "D17" = "0.04809"
    "D19" = "0.006001"
    "D21" = "0.0009828"
    "D23" = "3.715"
    "D24" = "2.230"
    "D25" = "0.3239"
    "D40" = "0.03873"
    "D41" = "0.006391"
    "D42" = "0.1073"
    "D43" = "0.002900"
    "D44" = "0.005072"
    "D45" = "0.00005626"
    "D47" = "0.9201"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) column (E) text updates.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"
